$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep their exact text representation
# (avoids Excel auto-converting numeric-looking strings to floats/sci notation).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.505.64'
$ws.Range("E2").Value = '  +4.04%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.735.29'
$ws.Range("E3").Value = '  +4.42%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.41'
$ws.Range("E5").Value = '  +3.61%  '

$ws.Range("E6").Value = '  +0.03%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4783'
$ws.Range("E7").Value = '  +3.56%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2658'
$ws.Range("E8").Value = '  +3.58%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06221'
$ws.Range("E9").Value = '  +1.44%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.735.98'
$ws.Range("E10").Value = '  +4.65%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07126'
$ws.Range("E11").Value = '  +2.73%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.70'
$ws.Range("E12").Value = '  +7.63%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6116'
$ws.Range("E13").Value = '  +7.54%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.523'
$ws.Range("E14").Value = '  +4.67%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '76.76'
$ws.Range("E15").Value = '  +2.33%  '

$ws.Range("E16").Value = '  +0.05%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.512.32'
$ws.Range("E17").Value = '  +4.06%  '

$ws.Range("E18").Value = '  +0.04%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000006889'
$ws.Range("E19").Value = '  +3.20%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.70'
$ws.Range("E20").Value = '  +3.09%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.959.74'
$ws.Range("E21").Value = '  +4.57%  '

$ws.Range("E22").Value = '  +3.58%  '

$ws.Range("E23").Value = '  +2.77%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.324'
$ws.Range("E24").Value = '  +2.32%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '135.85'
$ws.Range("E25").Value = '  +1.28%  '

$ws.Range("E26").Value = '  +3.13%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.791'
$ws.Range("E27").Value = '  +5.02%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.403'
$ws.Range("E28").Value = '  +2.68%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '106.36'
$ws.Range("E29").Value = '  +2.68%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.964'
$ws.Range("E30").Value = '  +0.76%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.701'
$ws.Range("E31").Value = '  +3.26%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.07886'
$ws.Range("E32").Value = '  +2.47%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04545'
$ws.Range("E33").Value = '  +4.98%  '

$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.615'
$ws.Range("E34").Value = '  +0.08%  '

$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6331'
$ws.Range("E35").Value = '  +5.78%  '

$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9888'
$ws.Range("E36").Value = '  +5.57%  '

$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9327'
$ws.Range("E37").Value = '  +2.44%  '

$ws.Range("B38").Value = 'Quant'
$ws.Range("C38").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '110.62'
$ws.Range("E38").Value = '  +2.83%  '

$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.981'
$ws.Range("E39").Value = '  +9.19%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.425'
$ws.Range("E40").Value = '  +1.71%  '

$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.005'
$ws.Range("E41").Value = '  +0.56%  '

$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.01508'
$ws.Range("E42").Value = '  +3.90%  '

$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.678'
$ws.Range("E43").Value = '  +13.87%  '

$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3892'
$ws.Range("E44").Value = '  +5.28%  '

$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.912'
$ws.Range("E45").Value = '  +13.43%  '

$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1188'
$ws.Range("E46").Value = '  +7.64%  '

$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05331'
$ws.Range("E47").Value = '  +1.36%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.927'
$ws.Range("E48").Value = '  +4.10%  '

$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '30.73'
$ws.Range("E49").Value = '  +1.32%  '

$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.255'
$ws.Range("E50").Value = '  +5.75%  '

$ws.Range("B51").Value = 'Decentraland'
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3420'
$ws.Range("E51").Value = '  +3.59%  '
